$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the extent of the sheet's data (21 columns x 55 rows = A1:U55)
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

# Rename header row cells: old format-version suffixes ("_old"/"_new") are
# replaced by the actual input-file names they represent ("_FV2310"/"_FV2404")
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# Freeze the header row (pane split after row 1) and select A2 like Excel does
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table (ListObject) named "Table1" so the
# renamed headers also become the table's column names, with filter buttons
$tblRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tblRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$wb.Save()
